# Build a small "database design" grid on Sheet1:
#   Column A = Caregories table, B = Products table, C = Customers table,
#   D = Orders table, E = Orders_Details table.
#
# Cells are filled in the exact order the original author typed them so that
# the shared-strings table comes out in the same order as the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: Caregories -------------------------------------------------
$ws.Range("A2").Value = "ID_cat"
$ws.Range("A3").Value = "Describtion_Cat"
$ws.Range("A1").Value = "Caregories"

# --- Column B: Products ----------------------------------------------------
$ws.Range("B1").Value = "Products"
$ws.Range("B2").Value = "ID_Product"
$ws.Range("B3").Value = "Label_Product"
$ws.Range("B4").Value = "QTE_in_Stock"
$ws.Range("B5").Value = "price )varchar)"
$ws.Range("B6").Value = "image_Product image"
$ws.Range("B7").Value = "ID_Cat"

# --- Column C: Customers -----------------------------------------------------
$ws.Range("C1").Value = "Customers"
$ws.Range("C2").Value = "ID_Customer"
$ws.Range("C3").Value = "First_Name"
$ws.Range("C4").Value = "Last_Name"
$ws.Range("C5").Value = "Tel (Ncahr(15)"
$ws.Range("C6").Value = "Email"
$ws.Range("C7").Value = "image_Customer nvarchar(max)"

# --- Column D: Orders --------------------------------------------------------
$ws.Range("D1").Value = "Orders"
$ws.Range("D2").Value = "ID_Order"
$ws.Range("D3").Value = "Date_Order DateTime"
$ws.Range("D4").Value = "Customer_ID"

# --- Column E: Orders_Details -------------------------------------------------
$ws.Range("E2").Value = "ID_Product"
$ws.Range("E3").Value = "ID_Order"
$ws.Range("E4").Value = "QTE"
$ws.Range("E1").Value = "Orders_Details"

# --- Column widths (closest achievable match to the source widths: 27.21875,
#     21.44140625 and 17.33203125 characters respectively) ---------------------
$ws.Columns.Item(1).ColumnWidth = 26.27
$ws.Columns.Item(2).ColumnWidth = 26.27
$ws.Columns.Item(3).ColumnWidth = 26.27
$ws.Columns.Item(4).ColumnWidth = 20.65
$ws.Columns.Item(5).ColumnWidth = 16.5

# --- Sheet view: selection on A2, scrolled so column C is left-most visible --
$null = $ws.Range("A2").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1

# --- Misc workbook window cosmetics (best effort) -----------------------------
$win.WindowState = -4140
